$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5910.7144
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 6062.5
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 6062.5
$ws.Range("M40").Value = -4825
$ws.Range("N40").Value = -6412.5
$ws.Range("H87").Value = 19999.953
$ws.Range("J87").Value = 19999.953
$ws.Range("L87").Value = 19999.953
$ws.Range("N87").Value = -22495.953
$ws.Range("H90").Value = 19999.953
$ws.Range("J90").Value = 19999.953
$ws.Range("L90").Value = 59999.859
$ws.Range("N90").Value = -72479.859
$ws.Range("H112").Value = 1594.8182
$ws.Range("J112").Value = 1597.037
$ws.Range("L112").Value = 4791.111
$ws.Range("N112").Value = -7007.111
$ws.Range("H138").Value = 2626.0334
$ws.Range("J138").Value = 4772.3076
$ws.Range("L138").Value = 14316.9228
$ws.Range("N138").Value = -24596.9228
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9123.574000000001
$ws.Range("I32").Value = 6768.8774
$ws.Range("K32").Value = 6768.8774
$ws.Range("M32").Value = -6481.8774
$ws.Range("H61").Value = 4290.64
$ws.Range("I61").Value = 2631.4375
$ws.Range("J61").Value = 7240.3335
$ws.Range("K61").Value = 2631.4375
$ws.Range("L61").Value = 7240.3335
$ws.Range("M61").Value = -2419.4375
$ws.Range("N61").Value = -7664.3335
$ws.Range("H74").Value = 3244.4546
$ws.Range("I74").Value = 3075.0557
$ws.Range("J74").Value = 4006.75
$ws.Range("K74").Value = 3075.0557
$ws.Range("L74").Value = 4006.75
$ws.Range("M74").Value = -2201.0557
$ws.Range("N74").Value = -5754.75
$ws.Range("H77").Value = 3244.4546
$ws.Range("I77").Value = 3075.0557
$ws.Range("J77").Value = 4006.75
$ws.Range("K77").Value = 15375.2785
$ws.Range("L77").Value = 20033.75
$ws.Range("M77").Value = -11007.2785
$ws.Range("N77").Value = -28769.75
$ws.Range("H88").Value = 2086.4167
$ws.Range("I88").Value = 2537.6667
$ws.Range("K88").Value = 2537.6667
$ws.Range("M88").Value = -2131.6667
$ws.Range("H91").Value = 2086.4167
$ws.Range("I91").Value = 2537.6667
$ws.Range("K91").Value = 2537.6667
$ws.Range("M91").Value = -1133.6667
$ws.Range("H114").Value = 60462
$ws.Range("J114").Value = 60462
$ws.Range("L114").Value = 60462
$ws.Range("N114").Value = -69140
$ws.Range("H132").Value = 34486620
$ws.Range("I132").Value = 62503404
$ws.Range("J132").Value = 4424.5386
$ws.Range("K132").Value = 187510212
$ws.Range("L132").Value = 13273.6158
$ws.Range("M132").Value = -187507682
$ws.Range("N132").Value = -18333.6158
$ws.Range("H134").Value = 59499.5
$ws.Range("J134").Value = 59499.5
$ws.Range("L134").Value = 59499.5
$ws.Range("N134").Value = -69639.5
$ws.Range("H135").Value = 79857.17999999999
$ws.Range("J135").Value = 79857.17999999999
$ws.Range("L135").Value = 79857.17999999999
$ws.Range("N135").Value = -89997.17999999999
$ws.Range("H136").Value = 4290.64
$ws.Range("I136").Value = 2631.4375
$ws.Range("J136").Value = 7240.3335
$ws.Range("K136").Value = 7894.3125
$ws.Range("L136").Value = 21721.0005
$ws.Range("M136").Value = -5344.3125
$ws.Range("N136").Value = -26821.0005
$ws.Range("H138").Value = 79998.164
$ws.Range("J138").Value = 79998.164
$ws.Range("L138").Value = 79998.164
$ws.Range("N138").Value = -90278.164
$ws.Range("H139").Value = 79857.336
$ws.Range("J139").Value = 79857.336
$ws.Range("L139").Value = 79857.336
$ws.Range("N139").Value = -90137.336
$ws.Range("H140").Value = 224935.48
$ws.Range("J140").Value = 224935.48
$ws.Range("L140").Value = 224935.48
$ws.Range("N140").Value = -235295.48
$ws.Range("H141").Value = 185927.84
$ws.Range("J141").Value = 185927.84
$ws.Range("L141").Value = 185927.84
$ws.Range("N141").Value = -196287.84

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H82").Value = 6848.2
$ws.Range("I82").Value = 6848.2
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 6848.2
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -6465.2
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 6848.2
$ws.Range("I85").Value = 6848.2
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 6848.2
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -5522.2
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 4020.7307
$ws.Range("I86").Value = 1569.9375
$ws.Range("J86").Value = 7942
$ws.Range("K86").Value = 1569.9375
$ws.Range("L86").Value = 7942
$ws.Range("M86").Value = -446.9375
$ws.Range("N86").Value = -10188
$ws.Range("H89").Value = 4020.7307
$ws.Range("I89").Value = 1569.9375
$ws.Range("J89").Value = 7942
$ws.Range("K89").Value = 7849.6875
$ws.Range("L89").Value = 39710
$ws.Range("M89").Value = -2233.6875
$ws.Range("N89").Value = -50942
$ws.Range("H94").Value = 2190.1428
$ws.Range("I94").Value = 2052.2632
$ws.Range("J94").Value = 3500
$ws.Range("K94").Value = 2052.2632
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = -1601.2632
$ws.Range("N94").Value = -4402
$ws.Range("H134").Value = 5423.5264
$ws.Range("I134").Value = 3748.111
$ws.Range("K134").Value = 11244.333
$ws.Range("M134").Value = -8709.332999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 19451224
$ws.Range("I99").Value = 16670752
$ws.Range("K99").Value = 16670752
$ws.Range("M99").Value = -16669254
$ws.Range("H126").Value = 19451224
$ws.Range("I126").Value = 16670752
$ws.Range("K126").Value = 50012256
$ws.Range("M126").Value = -50009786

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4965.467
$ws.Range("I132").Value = 3907.7144
$ws.Range("K132").Value = 11723.1432
$ws.Range("M132").Value = -9193.143199999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5618.8
$ws.Range("I7").Value = 2794.3333
$ws.Range("J7").Value = 6829.2856
$ws.Range("K7").Value = 2794.3333
$ws.Range("L7").Value = 6829.2856
$ws.Range("M7").Value = -2682.3333
$ws.Range("N7").Value = -7053.2856
$ws.Range("H40").Value = 11886.588
$ws.Range("I40").Value = 22222.8
$ws.Range("J40").Value = 7579.8335
$ws.Range("K40").Value = 22222.8
$ws.Range("L40").Value = 7579.8335
$ws.Range("M40").Value = -22086.8
$ws.Range("N40").Value = -7851.8335
$ws.Range("H93").Value = 373273.44
$ws.Range("I93").Value = 3098.5625
$ws.Range("J93").Value = 911709.6
$ws.Range("K93").Value = 3098.5625
$ws.Range("L93").Value = 911709.6
$ws.Range("M93").Value = -1850.5625
$ws.Range("N93").Value = -914205.6
$ws.Range("H126").Value = 5618.8
$ws.Range("I126").Value = 2794.3333
$ws.Range("J126").Value = 6829.2856
$ws.Range("K126").Value = 8382.999899999999
$ws.Range("L126").Value = 20487.8568
$ws.Range("M126").Value = -5912.999899999999
$ws.Range("N126").Value = -25427.8568
$ws.Range("H136").Value = 6476.44
$ws.Range("I136").Value = 4954
$ws.Range("J136").Value = 8414.091
$ws.Range("K136").Value = 14862
$ws.Range("L136").Value = 25242.273
$ws.Range("M136").Value = -12312
$ws.Range("N136").Value = -30342.273

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2386.4167
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 4700.1816
$ws.Range("I136").Value = 1560.4
$ws.Range("J136").Value = 7316.6665
$ws.Range("K136").Value = 4681.200000000001
$ws.Range("L136").Value = 21949.9995
$ws.Range("M136").Value = -2131.200000000001
$ws.Range("N136").Value = -27049.9995
